$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 95, pushing the former rows 95-128
# down to 97-130 (they keep their original data/formatting intact).
$ws.Rows.Item(95).Insert()
$ws.Rows.Item(95).Insert()

# Populate the new row 95 with the new weekly price entry.
$ws.Cells.Item(95,1).Value2  = 4
$ws.Cells.Item(95,2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(95,3).Value2  = "Los Lagos"
$ws.Cells.Item(95,4).Value2  = 44463
$ws.Cells.Item(95,5).Value2  = 10
$ws.Cells.Item(95,6).Value2  = "Fruta"
$ws.Cells.Item(95,7).Value2  = 100104
$ws.Cells.Item(95,8).Value2  = "Frutos de pepita"
$ws.Cells.Item(95,9).Value2  = 100104005
$ws.Cells.Item(95,10).Value2 = "Pera"
$ws.Cells.Item(95,11).Value2 = "Packham's Triumph"
$ws.Cells.Item(95,12).Value2 = "Primera"
$ws.Cells.Item(95,13).Value2 = 400
$ws.Cells.Item(95,14).Value2 = 16000
$ws.Cells.Item(95,15).Value2 = 16000
$ws.Cells.Item(95,16).Value2 = 16000
$ws.Cells.Item(95,17).Value2 = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(95,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(95,19).Value2 = 1067
$ws.Cells.Item(95,20).Value2 = 15

# Populate the new row 96 with the new weekly price entry.
$ws.Cells.Item(96,1).Value2  = 4
$ws.Cells.Item(96,2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(96,3).Value2  = "Los Lagos"
$ws.Cells.Item(96,4).Value2  = 44463
$ws.Cells.Item(96,5).Value2  = 10
$ws.Cells.Item(96,6).Value2  = "Fruta"
$ws.Cells.Item(96,7).Value2  = 100104
$ws.Cells.Item(96,8).Value2  = "Frutos de pepita"
$ws.Cells.Item(96,9).Value2  = 100104005
$ws.Cells.Item(96,10).Value2 = "Pera"
$ws.Cells.Item(96,11).Value2 = "Packham's Triumph"
$ws.Cells.Item(96,12).Value2 = "Segunda"
$ws.Cells.Item(96,13).Value2 = 200
$ws.Cells.Item(96,14).Value2 = 13000
$ws.Cells.Item(96,15).Value2 = 13000
$ws.Cells.Item(96,16).Value2 = 13000
$ws.Cells.Item(96,17).Value2 = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(96,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(96,19).Value2 = 867
$ws.Cells.Item(96,20).Value2 = 15
